# Grade workbook update:
#  - Student row 2's Final (G2) score corrected from 72 to 87.
#  - The curve formulas in columns I/J/K were reworked:
#      * I2 (the top scorer) now clamps the bonus baseline with an IF()
#        instead of a flat "+4.6" bump, so it can't get a bigger bonus than
#        everyone else just for being on top.
#      * I3:I25 use a flat "+1" bonus instead of "+4.6".
#      * J's "is the curve already maxed" threshold moved from 99.4 to 99.49.
#      * K's upper thresholds moved from 99.4 to 99.49, and the "bump to 60"
#        band moved from >=57.5 to >=53.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected Final score for the first student ---
$ws.Range("G2").Value = 87

# --- Column I: bonus-adjusted score ---
$ws.Range("I2").Formula = '=H2 + (IF($H$26 < 76.49, 76.49, $H$26) - $H$26)'
$ws.Range("I3:I25").Formula = '=H3 + (76.49 + 1 - $H$26)'

# --- Column J: curve so the top score hits 99.4 ---
$ws.Range("J2").Formula = '=IF(MAX($I$2:$I$25) < 99.49, I2 + (99.4 - MAX($I$2:$I$25)), I2)'
$ws.Range("J3:J25").Formula = '=IF(MAX($I$2:$I$25) < 99.49, I3 + (99.4 - MAX($I$2:$I$25)), I3)'

# --- Column K: final score with caps/floors applied ---
$ws.Range("K2").Formula = '=IF(AND(OR(F2=100,G2=100),J2>99.49),100,IF(AND(F2<100,G2<100,J2>99.49),99.49,IF(AND(J2>=53.5,J2<=59.4),60,IF(OR(F2=0,G2=0),0,J2))))'
$ws.Range("K3:K25").Formula = '=IF(AND(OR(F3=100,G3=100),J3>99.49),100,IF(AND(F3<100,G3<100,J3>99.49),99.49,IF(AND(J3>=53.5,J3<=59.4),60,IF(OR(F3=0,G3=0),0,J3))))'
